$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 21, shifting existing rows 21..126 down to 22..127
$ws.Rows.Item(21).Insert()

# Populate the newly inserted row 21 with the new weekly record
$ws.Range("A21").Value = 9
$ws.Range("B21").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C21").Value = "Metropolitana"
$ws.Range("D21").Value = 44453
$ws.Range("E21").Value = 13
$ws.Range("F21").Value = 100112026
$ws.Range("G21").Value = "Haba"
$ws.Range("H21").Value = "Sin especificar"
$ws.Range("I21").Value = "Primera"
$ws.Range("J21").Value = 52
$ws.Range("K21").Value = 15000
$ws.Range("L21").Value = 16000
$ws.Range("M21").Value = 15500
$ws.Range("N21").Value = '$/saco 25 kilos'
$ws.Range("O21").Value = "Provincia de Limarí"
$ws.Range("P21").Value = 620
$ws.Range("Q21").Value = 25
$ws.Range("R21").Value = "Hortaliza"
